$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[60.03741705920112, 67.0728691648524]"
$ws.Range("T2").Value = "[46.87178736618296, 51.59663851534458]"
$ws.Range("L3").Value = "[57.573043646227305, 69.02528619517061]"
$ws.Range("T3").Value = "[46.98972506106087, 54.00801576067052]"
